$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 168-169, pushing the existing rows 168-187 down to 170-189
$ws.Range("A168:A169").EntireRow.Insert()

# New row 168: Camote / 2a nueva(o) / Perú
$ws.Cells.Item(168, 1).Value = 4
$ws.Cells.Item(168, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(168, 3).Value = "Los Lagos"
$ws.Cells.Item(168, 4).Value = 44491
$ws.Cells.Item(168, 5).Value = 10
$ws.Cells.Item(168, 6).Value = 100112045
$ws.Cells.Item(168, 7).Value = "Zapallo"
$ws.Cells.Item(168, 8).Value = "Camote"
$ws.Cells.Item(168, 9).Value = "2a nueva(o)"
$ws.Cells.Item(168, 10).Value = 300
$ws.Cells.Item(168, 11).Value = 600
$ws.Cells.Item(168, 12).Value = 600
$ws.Cells.Item(168, 13).Value = 600
$ws.Cells.Item(168, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(168, 15).Value = "Perú"
$ws.Cells.Item(168, 16).Value = 600
$ws.Cells.Item(168, 17).Value = 1
$ws.Cells.Item(168, 18).Value = "Hortaliza"

# New row 169: Paine / 1a (guarda) / Región Metropolitana
$ws.Cells.Item(169, 1).Value = 4
$ws.Cells.Item(169, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(169, 3).Value = "Los Lagos"
$ws.Cells.Item(169, 4).Value = 44491
$ws.Cells.Item(169, 5).Value = 10
$ws.Cells.Item(169, 6).Value = 100112045
$ws.Cells.Item(169, 7).Value = "Zapallo"
$ws.Cells.Item(169, 8).Value = "Paine"
$ws.Cells.Item(169, 9).Value = "1a (guarda)"
$ws.Cells.Item(169, 10).Value = 800
$ws.Cells.Item(169, 11).Value = 400
$ws.Cells.Item(169, 12).Value = 400
$ws.Cells.Item(169, 13).Value = 400
$ws.Cells.Item(169, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(169, 15).Value = "Región Metropolitana"
$ws.Cells.Item(169, 16).Value = 400
$ws.Cells.Item(169, 17).Value = 1
$ws.Cells.Item(169, 18).Value = "Hortaliza"
